# Generate Report for Handoff
# Replace the old report GUID-named files/timestamps with the newly generated
# handoff report's identifiers across the Overview/zh-cn/de-de sheets, and
# keep the hyperlink display text in sync with the new cell text.

$wb = $excel.ActiveWorkbook

$oldGuidMd      = "ec1e5e54-8634-49c0-952d-08f711703ab5.md"
$newGuidMd      = "9dd8a24a-4dd4-4ca8-a56f-b6a3965886a3.md"

$oldZhXlf       = "ec1e5e54-8634-49c0-952d-08f711703ab5.318fa646bd55840d31c6c24a6a08f76d9e48541b.zh-cn.xlf"
$newZhXlf       = "9dd8a24a-4dd4-4ca8-a56f-b6a3965886a3.73519fdf91340b84ead8c1652dd0e2a0f2cdf09c.zh-cn.xlf"

$oldDeXlf       = "ec1e5e54-8634-49c0-952d-08f711703ab5.318fa646bd55840d31c6c24a6a08f76d9e48541b.de-de.xlf"
$newDeXlf       = "9dd8a24a-4dd4-4ca8-a56f-b6a3965886a3.73519fdf91340b84ead8c1652dd0e2a0f2cdf09c.de-de.xlf"

$oldZhDateTime  = "2016-03-10 19:01:57"
$newZhDateTime  = "2016-03-10 19:02:33"

$oldDeDateTime  = "2016-03-10 19:02:02"
$newDeDateTime  = "2016-03-10 19:02:37"

# --- Overview sheet: only the source markdown filename link ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuidMd

# --- zh-cn sheet: markdown filename, handoff xlf filename, handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newGuidMd
$wsZhCn.Range("C2").Value = $newZhXlf
$wsZhCn.Range("D2").Value = $newZhDateTime

# --- de-de sheet: markdown filename, handoff xlf filename, handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newGuidMd
$wsDeDe.Range("C2").Value = $newDeXlf
$wsDeDe.Range("D2").Value = $newDeDateTime

# --- Keep each hyperlink's displayed text in sync with the new cell text ---
foreach ($ws in $wb.Worksheets) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.TextToDisplay -eq $oldGuidMd) {
            $h.TextToDisplay = $newGuidMd
        } elseif ($h.TextToDisplay -eq $oldZhXlf) {
            $h.TextToDisplay = $newZhXlf
        } elseif ($h.TextToDisplay -eq $oldDeXlf) {
            $h.TextToDisplay = $newDeXlf
        }
    }
}
